$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REGULAR PICK")

# Header for the new UPH column
$ws.Range("C1").Value = "UPH"

# Rows 2-13 hold the per-user data; B currently stores negative quantities.
# Flip them to their absolute value and compute UPH = ABS(qty) * 60 / 169.
$factor = 60 / 169
for ($r = 2; $r -le 13; $r++) {
    $qty = [Math]::Abs($ws.Cells.Item($r, 2).Value2)
    $ws.Cells.Item($r, 2).Value = $qty
    $ws.Cells.Item($r, 3).Value = $qty * $factor
}
